$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5683
$ws.Range("L3").Value = 6180
$ws.Range("H4").Value = 1767
$ws.Range("L4").Value = 1522
$ws.Range("L5").Value = 366
$ws.Range("L6").Value = 5086
$ws.Range("H7").Value = 26083
$ws.Range("L7").Value = 18837

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 370
$ws.Range("L3").Value = 441
$ws.Range("L7").Value = 1249

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 303
$ws.Range("L6").Value = 243
$ws.Range("L7").Value = 857

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 85
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 215
$ws.Range("L3").Value = 253
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 721

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 115
$ws.Range("L7").Value = 361

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L3").Value = 13
$ws.Range("L5").Value = 67
$ws.Range("L7").Value = 607
$ws.Range("L8").Value = 1249
$ws.Range("L18").Value = 131
$ws.Range("L19").Value = 518
$ws.Range("L20").Value = 469
$ws.Range("L24").Value = 50
$ws.Range("L27").Value = 167
$ws.Range("L29").Value = 1063
$ws.Range("L33").Value = 857
$ws.Range("L36").Value = 240
$ws.Range("L37").Value = 721
$ws.Range("L41").Value = 82
$ws.Range("L42").Value = 609
$ws.Range("L43").Value = 137
$ws.Range("L44").Value = 132
$ws.Range("L45").Value = 36
$ws.Range("L47").Value = 123
$ws.Range("L48").Value = 245
$ws.Range("L49").Value = 97
$ws.Range("L53").Value = 203
$ws.Range("L54").Value = 416
$ws.Range("H63").Value = 317
$ws.Range("L63").Value = 56
$ws.Range("L64").Value = 122
$ws.Range("L65").Value = 361
$ws.Range("L66").Value = 55
$ws.Range("L67").Value = 653
$ws.Range("L72").Value = 79
$ws.Range("L74").Value = 15
$ws.Range("L76").Value = 291
$ws.Range("L77").Value = 126
$ws.Range("L79").Value = 516
$ws.Range("L85").Value = 933
$ws.Range("L89").Value = 263
$ws.Range("L90").Value = 195
$ws.Range("L94").Value = 233
$ws.Range("L95").Value = 267
$ws.Range("L96").Value = 214
$ws.Range("H101").Value = 26083
$ws.Range("L101").Value = 18837

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 253
$ws.Range("L7").Value = 653

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 202
$ws.Range("L7").Value = 416

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 314
$ws.Range("L6").Value = 267
$ws.Range("L7").Value = 1063

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 245

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 162
$ws.Range("L6").Value = 141
$ws.Range("L7").Value = 518

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 52
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 60
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 211
$ws.Range("L7").Value = 609

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 65
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 135
$ws.Range("L7").Value = 516

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L4").Value = 44
$ws.Range("L7").Value = 469

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L4").Value = 13
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 84
$ws.Range("L7").Value = 240

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 209
$ws.Range("L6").Value = 147
$ws.Range("L7").Value = 607

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 55
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L4").Value = 44
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 56
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L6").Value = 194
$ws.Range("L7").Value = 933

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 15
